$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 2
$ws.Range("A10").Value = 4
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 12

$ws.Range("F9:G10").FormulaArray = "=MMULT(A9:A10,C9:D9)"

$ws.Range("A12").Select()
